$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 12:39"

# --- Country name re-ordering (caused by re-sorting countries by case totals) ---
# Bielorrusia / Catar swap
$ws.Range("A30").Value = "Catar"
$ws.Range("A31").Value = "Bielorrusia"

# Finlandia / Marruecos swap
$ws.Range("A55").Value = "Marruecos"
$ws.Range("A56").Value = "Finlandia"

# Bolivia / Azerbaiyan / Uzbekistan / Camerun rotation
$ws.Range("A72").Value = "Camerun"
$ws.Range("A73").Value = "Bolivia"
$ws.Range("A74").Value = "Azerbaiyan"
$ws.Range("A75").Value = "Uzbekistan"

# --- Updated statistics values ---

# Row 5: Espana
$ws.Range("B5").Value = 264663
$ws.Range("C5").Value = 1880
$ws.Range("D5").Value = 176439
$ws.Range("E5").Value = 61603
$ws.Range("G5").Value = 143
$ws.Range("H5").Value = 26621

# Row 22: Suiza
$ws.Range("B22").Value = 30305
$ws.Range("C22").Value = 54
$ws.Range("E22").Value = 2075

# Row 30: Catar (new values)
$ws.Range("B30").Value = 22520
$ws.Range("C30").Value = 1189
$ws.Range("D30").Value = 2753
$ws.Range("E30").Value = 19753
$ws.Range("F30").Value = 72
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 14

# Row 31: Bielorrusia (new values)
$ws.Range("B31").Value = 22052
$ws.Range("D31").Value = 6050
$ws.Range("E31").Value = 15876
$ws.Range("F31").Value = 92
$ws.Range("H31").Value = 126

# Row 32: Emiratos Arabes Unidos
$ws.Range("B32").Value = 18198
$ws.Range("C32").Value = 781
$ws.Range("D32").Value = 4804
$ws.Range("E32").Value = 13196
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 198

# Row 55: Marruecos (new values)
$ws.Range("B55").Value = 6038
$ws.Range("C55").Value = 128
$ws.Range("D55").Value = 2545
$ws.Range("E55").Value = 3305
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 188

# Row 56: Finlandia (new values)
$ws.Range("B56").Value = 5963
$ws.Range("C56").Value = 83
$ws.Range("D56").Value = 4000
$ws.Range("E56").Value = 1696
$ws.Range("F56").Value = 45
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 267

# Row 72: Camerun (new values)
$ws.Range("B72").Value = 2579
$ws.Range("C72").Value = 305
$ws.Range("D72").Value = 1465
$ws.Range("E72").Value = 1000
$ws.Range("F72").Value = 28
$ws.Range("G72").Value = 6

# Row 73: Bolivia (new values)
$ws.Range("B73").Value = 2437
$ws.Range("C73").Value = 171
$ws.Range("D73").Value = 258
$ws.Range("E73").Value = 2065
$ws.Range("F73").Value = 3
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 114

# Row 74: Azerbaiyan (new values)
$ws.Range("B74").Value = 2422
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 1620
$ws.Range("E74").Value = 771
$ws.Range("F74").Value = 29
$ws.Range("H74").Value = 31

# Row 75: Uzbekistan (new values)
$ws.Range("B75").Value = 2411
$ws.Range("C75").Value = 62
$ws.Range("D75").Value = 1856
$ws.Range("E75").Value = 545
$ws.Range("F75").Value = 8
$ws.Range("H75").Value = 10

# Row 132: Montenegro
$ws.Range("E132").Value = 41
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 9
